# Add a "Correct_answer" column (D) to the sheet.
#   D1        = "Correct_answer"          (header)
#   D2:D17    = "l"                       (rows whose Color = Purple)
#   D18:D101  = "s"                       (rows whose Color = Blue)
# and move the active selection to D18:D101, scrolling the view so the
# new column is visible (mirrors the author's manual edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Correct_answer"
$ws.Range("D18:D101").Value = "s"
$ws.Range("D2:D17").Value = "l"

# Scroll/select like the original author did after adding the column.
$excel.Goto($ws.Range("A83"), $true)
$ws.Range("D18:D101").Select()
